$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2 values (recalibrated) ---
$ws.Range("B2").Value = 0.6337245365475435
$ws.Range("C2").Value = 0.07774148510829536
$ws.Range("D2").Value = 0.1718412239382025

# --- Append new monthly rows 266:277 (data through 2023) ---
$newData = @(
    @(44927, 1.147118638814245, 0.1139719076886744, 0.2307520696147135, 0.24903),
    @(44958, 1.155351166472892, 0.1189369534626325, 0.2358590252969582, 0.23755),
    @(44986, 1.159280358753098, 0.1235684903256264, 0.2621453884777878, 0.2074),
    @(45017, 1.163209551033304, 0.1231768212583258, 0.2553280272984077, 0.17092),
    @(45047, 1.166709031091589, 0.1216016047969423, 0.2555907188967245, 0.1307),
    @(45078, 1.170208511149873, 0.1222866409192133, 0.2669609706504602, 0.14978),
    @(45108, 1.173719182562309, 0.1229919207237912, 0.2509970053968685, 0.22925),
    @(45139, 1.177229853974744, 0.1304306391350569, 0.2545806997572367, 0.23604),
    @(45170, 1.179208384935429, 0.1401109389353736, 0.2544927629703319, 0.20583),
    @(45200, 1.181186915896113, 0.1405350854268179, 0.2690514056015423, 0.24936),
    @(45231, 1.179856082742833, 0.1251934609210341, 0.2503695190631033, 0.25939),
    @(45261, 1.178525249589552, 0.1157601160791554, 0.2470884693403191, 0.27367),
)

$startRow = 266
for ($i = 0; $i -lt $newData.Count; $i++) {
    $r = $startRow + $i
    $vals = $newData[$i]
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
    $ws.Cells.Item($r, 4).Value = $vals[3]
    $ws.Cells.Item($r, 5).Value = $vals[4]
}

# Replicate the date-column formatting (style) from the previous last row (A265)
$ws.Range("A265").Copy() | Out-Null
$ws.Range("A266:A277").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

Write-Host "Done"
